$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2066420664206642
$ws.Range("C2").Value = 0.5424354243542435
$ws.Range("J2").Value = 0.003690036900369004
$ws.Range("P2").Value = 0.1586715867158671
$ws.Range("S2").Value = 0.08856088560885608
$ws.Range("C3").Value = 0.01351351351351351
$ws.Range("J3").Value = 0.04054054054054054
$ws.Range("P3").Value = 0.7567567567567568
$ws.Range("S3").Value = 0.1891891891891892
$ws.Range("J4").Value = 0.02325581395348837
$ws.Range("O4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.7906976744186046
$ws.Range("S4").Value = 0.1627906976744186
$ws.Range("B6").Value = 0.08
$ws.Range("D6").Value = 0.01
$ws.Range("E6").Value = 0.005
$ws.Range("F6").Value = 0.075
$ws.Range("J6").Value = 0.27
$ws.Range("O6").Value = 0.035
$ws.Range("Q6").Value = 0.105
$ws.Range("R6").Value = 0.035
$ws.Range("S6").Value = 0.385
$ws.Range("B7").Value = 0.1015228426395939
$ws.Range("D7").Value = 0.01522842639593909
$ws.Range("F7").Value = 0.07106598984771574
$ws.Range("J7").Value = 0.1015228426395939
$ws.Range("O7").Value = 0.01015228426395939
$ws.Range("Q7").Value = 0.2030456852791878
$ws.Range("R7").Value = 0.07614213197969544
$ws.Range("S7").Value = 0.4213197969543147
$ws.Range("B8").Value = 0.09787234042553192
$ws.Range("D8").Value = 0.01702127659574468
$ws.Range("E8").Value = 0.00425531914893617
$ws.Range("F8").Value = 0.05531914893617021
$ws.Range("J8").Value = 0.1063829787234043
$ws.Range("O8").Value = 0.01702127659574468
$ws.Range("Q8").Value = 0.1808510638297872
$ws.Range("R8").Value = 0.09574468085106383
$ws.Range("S8").Value = 0.425531914893617
$ws.Range("B9").Value = 0.1061946902654867
$ws.Range("D9").Value = 0.008849557522123894
$ws.Range("F9").Value = 0.04867256637168142
$ws.Range("J9").Value = 0.1017699115044248
$ws.Range("O9").Value = 0.03097345132743363
$ws.Range("Q9").Value = 0.2256637168141593
$ws.Range("R9").Value = 0.084070796460177
$ws.Range("S9").Value = 0.3938053097345133
$ws.Range("B10").Value = 0.09200343938091143
$ws.Range("D10").Value = 0.02407566638005159
$ws.Range("E10").Value = 0.001719690455717971
$ws.Range("F10").Value = 0.05846947549441101
$ws.Range("J10").Value = 0.1358555460017197
$ws.Range("O10").Value = 0.02235597592433362
$ws.Range("Q10").Value = 0.2252794496990542
$ws.Range("R10").Value = 0.07222699914015478
$ws.Range("S10").Value = 0.3680137575236457
$ws.Range("G11").Value = 0.09881422924901186
$ws.Range("J11").Value = 0.05928853754940711
$ws.Range("K11").Value = 0.1383399209486166
$ws.Range("L11").Value = 0.6877470355731226
$ws.Range("S11").Value = 0.0158102766798419
$ws.Range("G12").Value = 0.8248587570621468
$ws.Range("J12").Value = 0.1468926553672316
$ws.Range("L12").Value = 0.02259887005649718
$ws.Range("S12").Value = 0.005649717514124294
$ws.Range("G13").Value = 0.7872340425531915
$ws.Range("J13").Value = 0.1914893617021277
$ws.Range("S13").Value = 0.02127659574468085
$ws.Range("F15").Value = 0.0199203187250996
$ws.Range("H15").Value = 0.1872509960159363
$ws.Range("I15").Value = 0.07171314741035857
$ws.Range("J15").Value = 0.3147410358565737
$ws.Range("K15").Value = 0.06772908366533864
$ws.Range("M15").Value = 0.01593625498007968
$ws.Range("O15").Value = 0.07569721115537849
$ws.Range("S15").Value = 0.2470119521912351
$ws.Range("F16").Value = 0.0160427807486631
$ws.Range("H16").Value = 0.2513368983957219
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.3368983957219251
$ws.Range("K16").Value = 0.1176470588235294
$ws.Range("M16").Value = 0.03208556149732621
$ws.Range("O16").Value = 0.06951871657754011
$ws.Range("S16").Value = 0.0855614973262032
$ws.Range("F17").Value = 0.02863436123348018
$ws.Range("H17").Value = 0.1938325991189427
$ws.Range("I17").Value = 0.1255506607929515
$ws.Range("J17").Value = 0.3722466960352423
$ws.Range("K17").Value = 0.09251101321585903
$ws.Range("M17").Value = 0.013215859030837
$ws.Range("O17").Value = 0.06828193832599119
$ws.Range("S17").Value = 0.105726872246696
$ws.Range("F18").Value = 0.01785714285714286
$ws.Range("H18").Value = 0.1785714285714286
$ws.Range("I18").Value = 0.08333333333333333
$ws.Range("J18").Value = 0.4404761904761905
$ws.Range("K18").Value = 0.06547619047619048
$ws.Range("M18").Value = 0.01785714285714286
$ws.Range("O18").Value = 0.1130952380952381
$ws.Range("S18").Value = 0.08333333333333333
$ws.Range("F19").Value = 0.01357082273112807
$ws.Range("H19").Value = 0.2196776929601357
$ws.Range("I19").Value = 0.1009329940627651
$ws.Range("J19").Value = 0.368108566581849
$ws.Range("K19").Value = 0.1043256997455471
$ws.Range("M19").Value = 0.02374893977947413
$ws.Range("O19").Value = 0.07888040712468193
$ws.Range("S19").Value = 0.09075487701441901
